$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 3 rows (old rows 8-10, MuSCs as sending cluster) entirely
$ws.Rows("8:10").Delete() | Out-Null

# Update remaining data rows (2-7) with the new TPM-based values
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.919643
$ws.Range("H2").Value = 53.75892899999999
$ws.Range("I2").Value = 0.8982899767221961
$ws.Range("J2").Value = 0.8982899767221962
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.097779
$ws.Range("N2").Value = 78.29333700000001
$ws.Range("O2").Value = 0.9922055808976035
$ws.Range("P2").Value = 0.9922055808976036
$ws.Range("Q2").Value = 467.662882772897
$ws.Range("R2").Value = 4208.965944956073
$ws.Range("S2").Value = 0.8912883281681413
$ws.Range("T2").Value = 0.8912883281681415

$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.919643
$ws.Range("H3").Value = 53.75892899999999
$ws.Range("I3").Value = 0.8982899767221961
$ws.Range("J3").Value = 0.8982899767221962
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.205015
$ws.Range("N3").Value = 0.6150450000000001
$ws.Range("O3").Value = 0.007794419102396499
$ws.Range("P3").Value = 0.007794419102396499
$ws.Range("Q3").Value = 3.673795609645
$ws.Range("R3").Value = 33.064160486805
$ws.Range("S3").Value = 0.007001648554054792
$ws.Range("T3").Value = 0.007001648554054793

$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.359006333333333
$ws.Range("H4").Value = 4.077019
$ws.Range("I4").Value = 0.06812533974785755
$ws.Range("J4").Value = 0.06812533974785755
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.097779
$ws.Range("N4").Value = 78.29333700000001
$ws.Range("O4").Value = 0.9922055808976035
$ws.Range("P4").Value = 0.9922055808976036
$ws.Range("Q4").Value = 35.46704694693367
$ws.Range("R4").Value = 319.203422522403
$ws.Range("S4").Value = 0.06759434229836959
$ws.Range("T4").Value = 0.06759434229836959

$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.359006333333333
$ws.Range("H5").Value = 4.077019
$ws.Range("I5").Value = 0.06812533974785755
$ws.Range("J5").Value = 0.06812533974785755
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.205015
$ws.Range("N5").Value = 0.6150450000000001
$ws.Range("O5").Value = 0.007794419102396499
$ws.Range("P5").Value = 0.007794419102396499
$ws.Range("Q5").Value = 0.2786166834283333
$ws.Range("R5").Value = 2.507550150855
$ws.Range("S5").Value = 0.0005309974494879524
$ws.Range("T5").Value = 0.0005309974494879524

$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.669968
$ws.Range("H6").Value = 2.009904
$ws.Range("I6").Value = 0.03358468352994624
$ws.Range("J6").Value = 0.03358468352994624
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 26.097779
$ws.Range("N6").Value = 78.29333700000001
$ws.Range("O6").Value = 0.9922055808976035
$ws.Range("P6").Value = 0.9922055808976036
$ws.Range("Q6").Value = 17.484676801072
$ws.Range("R6").Value = 157.362091209648
$ws.Range("S6").Value = 0.03332291043109248
$ws.Range("T6").Value = 0.03332291043109249

$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.669968
$ws.Range("H7").Value = 2.009904
$ws.Range("I7").Value = 0.03358468352994624
$ws.Range("J7").Value = 0.03358468352994624
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.205015
$ws.Range("N7").Value = 0.6150450000000001
$ws.Range("O7").Value = 0.007794419102396499
$ws.Range("P7").Value = 0.007794419102396499
$ws.Range("Q7").Value = 0.13735348952
$ws.Range("R7").Value = 1.23618140568
$ws.Range("S7").Value = 0.000261773098853754
$ws.Range("T7").Value = 0.0002617730988537541
